$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy cell formatting (style) down from the last existing data row (90) ---
# Column A (group) and Column B (name) for every new row 91:144
$ws.Range("A90").Copy()
$ws.Range("A91:A144").PasteSpecial(-4122)
$ws.Range("B90").Copy()
$ws.Range("B91:B144").PasteSpecial(-4122)

# Column C (comment) cells that need values: copy format from an existing C-styled cell
$ws.Range("C84").Copy()
$ws.Range("C103").PasteSpecial(-4122)

# Column D (color) cells that need values: copy format from an existing D-styled cell
$ws.Range("D71").Copy()
$ws.Range("D101").PasteSpecial(-4122)
$ws.Range("D71").Copy()
$ws.Range("D108").PasteSpecial(-4122)
$ws.Range("D71").Copy()
$ws.Range("D140").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Assign cell values row-by-row, left-to-right, to match natural entry order ---
$ws.Cells.Item(91, 1).Value = 2
$ws.Cells.Item(91, 2).Value = "你的薛爸爸"
$ws.Cells.Item(92, 1).Value = 2
$ws.Cells.Item(92, 2).Value = "妖怪蜀黍"
$ws.Cells.Item(93, 1).Value = 2
$ws.Cells.Item(93, 2).Value = "刘沈钧"
$ws.Cells.Item(94, 1).Value = 2
$ws.Cells.Item(94, 2).Value = "饭饭"
$ws.Cells.Item(95, 1).Value = 2
$ws.Cells.Item(95, 2).Value = "加利斯"
$ws.Cells.Item(96, 1).Value = 2
$ws.Cells.Item(96, 2).Value = "滑稽拌饭"
$ws.Cells.Item(97, 1).Value = 2
$ws.Cells.Item(97, 2).Value = "敖鹏柽"
$ws.Cells.Item(98, 1).Value = 2
$ws.Cells.Item(98, 2).Value = "张逸尘"
$ws.Cells.Item(99, 1).Value = 2
$ws.Cells.Item(99, 2).Value = "红茶坊"
$ws.Cells.Item(100, 1).Value = 2
$ws.Cells.Item(100, 2).Value = "小灰"
$ws.Cells.Item(101, 1).Value = 2
$ws.Cells.Item(101, 2).Value = "七海Nana7mi"
$ws.Cells.Item(101, 4).Value = "#ff70cb"
$ws.Cells.Item(102, 1).Value = 2
$ws.Cells.Item(102, 2).Value = "松松小笨蛋"
$ws.Cells.Item(103, 1).Value = 1
$ws.Cells.Item(103, 2).Value = "NeroCloud"
$ws.Cells.Item(103, 3).Value = "挺好玩的小游戏"
$ws.Cells.Item(104, 1).Value = 2
$ws.Cells.Item(104, 2).Value = "策划大大的大跌"
$ws.Cells.Item(105, 1).Value = 2
$ws.Cells.Item(105, 2).Value = "tiny cherry petal"
$ws.Cells.Item(106, 1).Value = 2
$ws.Cells.Item(106, 2).Value = "啊哈哈"
$ws.Cells.Item(107, 1).Value = 2
$ws.Cells.Item(107, 2).Value = "多多"
$ws.Cells.Item(108, 1).Value = 2
$ws.Cells.Item(108, 2).Value = "般若鬼面 "
$ws.Cells.Item(108, 4).Value = "#921AFF"
$ws.Cells.Item(109, 1).Value = 2
$ws.Cells.Item(109, 2).Value = "璟年"
$ws.Cells.Item(110, 1).Value = 2
$ws.Cells.Item(110, 2).Value = "FAY"
$ws.Cells.Item(111, 1).Value = 2
$ws.Cells.Item(111, 2).Value = "干死黄旭东"
$ws.Cells.Item(112, 1).Value = 2
$ws.Cells.Item(112, 2).Value = "五胖是真的"
$ws.Cells.Item(113, 1).Value = 2
$ws.Cells.Item(113, 2).Value = "狗贼给箱子"
$ws.Cells.Item(114, 1).Value = 2
$ws.Cells.Item(114, 2).Value = "尹东洋"
$ws.Cells.Item(115, 1).Value = 2
$ws.Cells.Item(115, 2).Value = "红阳"
$ws.Cells.Item(116, 1).Value = 2
$ws.Cells.Item(116, 2).Value = "barusamikosu"
$ws.Cells.Item(117, 1).Value = 2
$ws.Cells.Item(117, 2).Value = "小保安"
$ws.Cells.Item(118, 1).Value = 2
$ws.Cells.Item(118, 2).Value = "Tony"
$ws.Cells.Item(119, 1).Value = 2
$ws.Cells.Item(119, 2).Value = "戒（ ）"
$ws.Cells.Item(120, 1).Value = 2
$ws.Cells.Item(120, 2).Value = "下半鸭"
$ws.Cells.Item(121, 1).Value = 2
$ws.Cells.Item(121, 2).Value = "Xiaosaye"
$ws.Cells.Item(122, 1).Value = 2
$ws.Cells.Item(122, 2).Value = "NAKAS"
$ws.Cells.Item(123, 1).Value = 2
$ws.Cells.Item(123, 2).Value = "icefrog"
$ws.Cells.Item(124, 1).Value = 2
$ws.Cells.Item(124, 2).Value = "南怪"
$ws.Cells.Item(125, 1).Value = 2
$ws.Cells.Item(125, 2).Value = "梁一斤"
$ws.Cells.Item(126, 1).Value = 2
$ws.Cells.Item(126, 2).Value = "熊公子nice哥"
$ws.Cells.Item(127, 1).Value = 2
$ws.Cells.Item(127, 2).Value = "徙南"
$ws.Cells.Item(128, 1).Value = 2
$ws.Cells.Item(128, 2).Value = "粟粟"
$ws.Cells.Item(129, 1).Value = 2
$ws.Cells.Item(129, 2).Value = "王哲"
$ws.Cells.Item(130, 1).Value = 2
$ws.Cells.Item(130, 2).Value = "HPENG"
$ws.Cells.Item(131, 1).Value = 2
$ws.Cells.Item(131, 2).Value = "残梦空恨"
$ws.Cells.Item(132, 1).Value = 2
$ws.Cells.Item(132, 2).Value = "测试名字的长度"
$ws.Cells.Item(133, 1).Value = 2
$ws.Cells.Item(133, 2).Value = "夏吉野"
$ws.Cells.Item(134, 1).Value = 2
$ws.Cells.Item(134, 2).Value = "你的名字"
$ws.Cells.Item(135, 1).Value = 2
$ws.Cells.Item(135, 2).Value = "幻想丶当时"
$ws.Cells.Item(136, 1).Value = 2
$ws.Cells.Item(136, 2).Value = "琦琦"
$ws.Cells.Item(137, 1).Value = 2
$ws.Cells.Item(137, 2).Value = "执笔诉情"
$ws.Cells.Item(138, 1).Value = 2
$ws.Cells.Item(138, 2).Value = "陈声"
$ws.Cells.Item(139, 1).Value = 2
$ws.Cells.Item(139, 2).Value = "恶魔在右丷"
$ws.Cells.Item(140, 1).Value = 2
$ws.Cells.Item(140, 2).Value = "苏苏"
$ws.Cells.Item(140, 4).Value = "pink"
$ws.Cells.Item(141, 1).Value = 2
$ws.Cells.Item(141, 2).Value = "哈哈和大叔"
$ws.Cells.Item(142, 1).Value = 2
$ws.Cells.Item(142, 2).Value = "攀大大"
$ws.Cells.Item(143, 1).Value = 2
$ws.Cells.Item(143, 2).Value = "电饭锅"
$ws.Cells.Item(144, 1).Value = 2
$ws.Cells.Item(144, 2).Value = "梦鼠"

# --- Selection / view state ---
$ws.Range("C144").Select()
